# Update "想去人数" (interested-people count) values in the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 445
$ws1.Range("F8").Value = 14025
$ws1.Range("F9").Value = 109
$ws1.Range("F11").Value = 5650
$ws1.Range("F15").Value = 51
$ws1.Range("F16").Value = 1225
$ws1.Range("F19").Value = 761
$ws1.Range("F20").Value = 2911
$ws1.Range("F21").Value = 47
$ws1.Range("F22").Value = 10431
$ws1.Range("F25").Value = 53
$ws1.Range("F26").Value = 3709

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 445
$ws4.Range("F9").Value = 14025
$ws4.Range("F10").Value = 109
$ws4.Range("F12").Value = 5650
$ws4.Range("F16").Value = 51
$ws4.Range("F17").Value = 1225
$ws4.Range("F20").Value = 761
$ws4.Range("F21").Value = 2911
$ws4.Range("F22").Value = 47
$ws4.Range("F24").Value = 10431
$ws4.Range("F27").Value = 53
$ws4.Range("F28").Value = 3709
